$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "61.621.88"
$c.Style = $s

$c = $ws.Range("E2")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.38%  "
$c.Style = $s

$c = $ws.Range("D3")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "3.452.09"
$c.Style = $s

$c = $ws.Range("E3")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.76%  "
$c.Style = $s

$c = $ws.Range("D4")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = $s

$c = $ws.Range("E4")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.01%  "
$c.Style = $s

$c = $ws.Range("D5")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "579.70"
$c.Style = $s

$c = $ws.Range("E5")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.34%  "
$c.Style = $s

$c = $ws.Range("D6")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "147.66"
$c.Style = $s

$c = $ws.Range("E6")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +8.02%  "
$c.Style = $s

$c = $ws.Range("D7")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "3.453.76"
$c.Style = $s

$c = $ws.Range("E7")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.84%  "
$c.Style = $s

$c = $ws.Range("E8")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.07%  "
$c.Style = $s

$c = $ws.Range("D10")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "7.67"
$c.Style = $s

$c = $ws.Range("E10")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.60%  "
$c.Style = $s

$c = $ws.Range("E11")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -0.02%  "
$c.Style = $s

$c = $ws.Range("E12")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -0.31%  "
$c.Style = $s

$c = $ws.Range("D13")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "4.045.34"
$c.Style = $s

$c = $ws.Range("E13")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.90%  "
$c.Style = $s

$c = $ws.Range("D14")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "27.84"
$c.Style = $s

$c = $ws.Range("E14")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +7.29%  "
$c.Style = $s

$c = $ws.Range("E15")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -0.86%  "
$c.Style = $s

$c = $ws.Range("E16")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.15%  "
$c.Style = $s

$c = $ws.Range("D17")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "3.448.54"
$c.Style = $s

$c = $ws.Range("E17")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.60%  "
$c.Style = $s

$c = $ws.Range("D18")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "61.692.24"
$c.Style = $s

$c = $ws.Range("E18")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.41%  "
$c.Style = $s

$c = $ws.Range("D19")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "6.32"
$c.Style = $s

$c = $ws.Range("E19")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +7.78%  "
$c.Style = $s

$c = $ws.Range("E20")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.00%  "
$c.Style = $s

$c = $ws.Range("E21")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.71%  "
$c.Style = $s

$c = $ws.Range("D22")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "385.01"
$c.Style = $s

$c = $ws.Range("E22")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.46%  "
$c.Style = $s

$c = $ws.Range("D23")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.566"
$c.Style = $s

$c = $ws.Range("E23")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.58%  "
$c.Style = $s

$c = $ws.Range("D24")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "3.588.51"
$c.Style = $s

$c = $ws.Range("E24")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.49%  "
$c.Style = $s

$c = $ws.Range("E25")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -0.06%  "
$c.Style = $s

$c = $ws.Range("E26")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.49%  "
$c.Style = $s

$c = $ws.Range("E27")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.68%  "
$c.Style = $s

$c = $ws.Range("E28")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -0.92%  "
$c.Style = $s

$c = $ws.Range("E29")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +7.48%  "
$c.Style = $s

$c = $ws.Range("D30")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "7.80"
$c.Style = $s

$c = $ws.Range("E30")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +4.08%  "
$c.Style = $s

$c = $ws.Range("D31")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.56"
$c.Style = $s

$c = $ws.Range("E31")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -13.99%  "
$c.Style = $s

$c = $ws.Range("E32")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -0.01%  "
$c.Style = $s

$c = $ws.Range("E33")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.48%  "
$c.Style = $s

$c = $ws.Range("D34")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.18"
$c.Style = $s

$c = $ws.Range("E34")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.28%  "
$c.Style = $s

$c = $ws.Range("E35")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -0.03%  "
$c.Style = $s

$c = $ws.Range("D36")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "24.03"
$c.Style = $s

$c = $ws.Range("E36")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.85%  "
$c.Style = $s

$c = $ws.Range("D37")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "5.24"
$c.Style = $s

$c = $ws.Range("E37")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.87%  "
$c.Style = $s

$c = $ws.Range("E38")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +3.39%  "
$c.Style = $s

$c = $ws.Range("E39")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.14%  "
$c.Style = $s

$c = $ws.Range("D40")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "166.03"
$c.Style = $s

$c = $ws.Range("E40")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.61%  "
$c.Style = $s

$c = $ws.Range("D41")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0786"
$c.Style = $s

$c = $ws.Range("E41")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.81%  "
$c.Style = $s

$c = $ws.Range("D42")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "26.25"
$c.Style = $s

$c = $ws.Range("E42")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +9.41%  "
$c.Style = $s

$c = $ws.Range("D43")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.792"
$c.Style = $s

$c = $ws.Range("E43")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.69%  "
$c.Style = $s

$c = $ws.Range("D44")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = $s

$c = $ws.Range("E44")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.02%  "
$c.Style = $s

$c = $ws.Range("E45")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.72%  "
$c.Style = $s

$c = $ws.Range("E46")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.37%  "
$c.Style = $s

$c = $ws.Range("D47")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "41.99"
$c.Style = $s

$c = $ws.Range("E47")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.36%  "
$c.Style = $s

$c = $ws.Range("D48")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.620.89"
$c.Style = $s

$c = $ws.Range("E48")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +11.25%  "
$c.Style = $s

$c = $ws.Range("E49")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -2.32%  "
$c.Style = $s

$c = $ws.Range("D50")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "23.73"
$c.Style = $s

$c = $ws.Range("E50")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +3.95%  "
$c.Style = $s

$c = $ws.Range("E51")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.32%  "
$c.Style = $s

